# Regenerate s_val data to filter save games.
# Updates the computed TB/d2S/K/IP/sum values for each row (2-6) on Sheet1.
# The "Win" column (F) and date column (A) are left unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (2023-08-21)
$ws.Range("B2").Value = 0.7287194209349384
$ws.Range("C2").Value = 0.3375848360084654
$ws.Range("D2").Value = 16.98373111632243
$ws.Range("E2").Value = 0.4998867070740569
$ws.Range("G2").Value = 18.54992208033989

# Row 3 (2023-05-30)
$ws.Range("B3").Value = 0.001754667048134761
$ws.Range("C3").Value = 0.004309184025731883
$ws.Range("D3").Value = 116886.6739907443
$ws.Range("E3").Value = 246.9852506941017
$ws.Range("G3").Value = 117133.6653052895

# Row 4 (2023-04-17)
$ws.Range("B4").Value = 1.505614041169197
$ws.Range("C4").Value = 0.3375848360084654
$ws.Range("D4").Value = 0.7127328510149897
$ws.Range("E4").Value = 0.4998867070740569
$ws.Range("G4").Value = 3.055818435266709

# Row 5 (2023-03-10)
$ws.Range("B5").Value = 0.1554434735375247
$ws.Range("C5").Value = 0.3375848360084654
$ws.Range("D5").Value = 3.082599426703578
$ws.Range("E5").Value = 246.9852506941017
$ws.Range("G5").Value = 250.5608784303512

# Row 6 (2023-03-03)
$ws.Range("B6").Value = 3.182878228561681
$ws.Range("C6").Value = 1.65323645889881
$ws.Range("D6").Value = 16.98373111632243
$ws.Range("E6").Value = 6.48142807727062
$ws.Range("G6").Value = 28.30127388105354
